# Visual updates to the app
# - Rename the "Cost" category label to "Expenditure" throughout the payment sheets.
# - Fix a mis-tagged grocery row (S-market, row 22) that was recorded as "Income".
# - Correct the October salary figure on the OP sheet (3000 -> 2000).
# - Refresh the remembered cell selection on the OP and Nordea sheets.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("OP", "Nordea")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ($cell.Value2 -eq "Cost") {
                $cell.Value = "Expenditure"
            }
        }
    }
}

$wsOP = $wb.Worksheets.Item("OP")

# Row 22 (grocery purchase at S-market) was wrongly marked as Income
$wsOP.Range("C22").Value = "Expenditure"

# October salary payment amount correction
$wsOP.Range("D2").Value = 2000

# Restore the active-cell selections recorded in each sheet's view
[void]$wsOP.Range("G15").Select()

$wsNordea = $wb.Worksheets.Item("Nordea")
[void]$wsNordea.Range("C11").Select()

[void]$wb.Worksheets.Item("OP").Activate()
